# Apply updated crypto price/volume values (Mon Oct 16 20:47:29 UTC 2023 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.517.28'
$ws.Range("E2").Value = '  +4.40%  '
$ws.Range("D3").Value = '1.592.63'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.78'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.96'
$ws.Range("E8").Value = '  +8.17%  '
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").Value = '1.820.33'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").Value = '1.588.21'
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.81'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").Value = '28.515.15'
$ws.Range("E16").Value = '  +4.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.02'
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.61'
$ws.Range("E18").Value = '  +7.17%  '
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.42'
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.04'
$ws.Range("E25").Value = '  -1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.32'
$ws.Range("E26").Value = '  +1.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.63'
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.16'
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("D34").Value = '1.413.35'
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("E36").Value = '  -6.28%  '
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.56'
$ws.Range("E39").Value = '  +6.18%  '
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.76'
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.82'
$ws.Range("E45").Value = '  +5.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.76'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '1.733.64'
$ws.Range("E47").Value = '  +1.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.90'
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  +4.39%  '
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.68'
$ws.Range("E51").Value = '  +16.75%  '
